$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "76.273.13"
Set-TextValue $ws.Range("E2") "  -0.50%  "
Set-TextValue $ws.Range("D3") "3.080.38"
Set-TextValue $ws.Range("E3") "  +4.43%  "
Set-TextValue $ws.Range("E4") "  -0.10%  "
Set-TextValue $ws.Range("D5") "197.88"
Set-TextValue $ws.Range("E5") "  -0.61%  "
Set-TextValue $ws.Range("D6") "614.92"
Set-TextValue $ws.Range("E6") "  +3.27%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.551"
Set-TextValue $ws.Range("E8") "  +0.13%  "
Set-TextValue $ws.Range("D9") "0.210"
Set-TextValue $ws.Range("E9") "  +7.07%  "
Set-TextValue $ws.Range("D10") "3.078.29"
Set-TextValue $ws.Range("E10") "  +4.24%  "
Set-TextValue $ws.Range("D11") "0.441"
Set-TextValue $ws.Range("E11") "  -0.29%  "
Set-TextValue $ws.Range("E12") "  -0.12%  "
Set-TextValue $ws.Range("D13") "5.22"
Set-TextValue $ws.Range("E13") "  +6.83%  "
Set-TextValue $ws.Range("D14") "3.646.63"
Set-TextValue $ws.Range("E14") "  +4.31%  "
Set-TextValue $ws.Range("D15") "29.19"
Set-TextValue $ws.Range("E15") "  +3.15%  "
Set-TextValue $ws.Range("D16") "76.131.98"
Set-TextValue $ws.Range("E16") "  -0.65%  "
Set-TextValue $ws.Range("D17") "0.0000194"
Set-TextValue $ws.Range("E17") "  +2.82%  "
Set-TextValue $ws.Range("D18") "3.079.13"
Set-TextValue $ws.Range("E18") "  +3.97%  "
Set-TextValue $ws.Range("D19") "13.58"
Set-TextValue $ws.Range("E19") "  +0.11%  "
Set-TextValue $ws.Range("D20") "9.13"
Set-TextValue $ws.Range("E20") "  +5.04%  "
Set-TextValue $ws.Range("D21") "380.71"
Set-TextValue $ws.Range("E21") "  +1.65%  "
Set-TextValue $ws.Range("D22") "2.54"
Set-TextValue $ws.Range("E22") "  +12.68%  "
Set-TextValue $ws.Range("D23") "4.42"
Set-TextValue $ws.Range("E23") "  +2.31%  "
Set-TextValue $ws.Range("D24") "3.239.49"
Set-TextValue $ws.Range("E24") "  +4.29%  "
Set-TextValue $ws.Range("E25") "  -0.33%  "
Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +0.10%  "
Set-TextValue $ws.Range("D27") "4.36"
Set-TextValue $ws.Range("E27") "  +2.00%  "
Set-TextValue $ws.Range("D28") "9.94"
Set-TextValue $ws.Range("E28") "  +2.96%  "
Set-TextValue $ws.Range("D29") "0.0000108"
Set-TextValue $ws.Range("E29") "  +1.07%  "
Set-TextValue $ws.Range("D30") "0.997"
Set-TextValue $ws.Range("E30") "  -0.49%  "
Set-TextValue $ws.Range("D31") "8.32"
Set-TextValue $ws.Range("E31") "  -0.03%  "
Set-TextValue $ws.Range("D32") "1.42"
Set-TextValue $ws.Range("E32") "  +2.80%  "
Set-TextValue $ws.Range("D33") "497.82"
Set-TextValue $ws.Range("E33") "  -0.31%  "
Set-TextValue $ws.Range("E34") "  +4.56%  "
Set-TextValue $ws.Range("E35") "  -0.13%  "
Set-TextValue $ws.Range("D36") "20.79"
Set-TextValue $ws.Range("E36") "  +2.95%  "
Set-TextValue $ws.Range("D37") "0.123"
Set-TextValue $ws.Range("E37") "  +11.64%  "
Set-TextValue $ws.Range("D38") "161.95"
Set-TextValue $ws.Range("E38") "  -1.87%  "
Set-TextValue $ws.Range("D39") "194.93"
Set-TextValue $ws.Range("E39") "  +7.99%  "
Set-TextValue $ws.Range("D40") "20.07"
Set-TextValue $ws.Range("E40") "  +0.63%  "
Set-TextValue $ws.Range("D41") "0.378"
Set-TextValue $ws.Range("E41") "  -3.97%  "
Set-TextValue $ws.Range("D42") "0.103"
Set-TextValue $ws.Range("E42") "  -8.75%  "
Set-TextValue $ws.Range("E43") "  +0.03%  "
Set-TextValue $ws.Range("D44") "0.800"
Set-TextValue $ws.Range("E44") "  +21.51%  "
Set-TextValue $ws.Range("D45") "5.11"
Set-TextValue $ws.Range("E45") "  +3.79%  "
Set-TextValue $ws.Range("D46") "1.25"
Set-TextValue $ws.Range("E46") "  +5.30%  "
Set-TextValue $ws.Range("D47") "41.31"
Set-TextValue $ws.Range("E47") "  +2.91%  "
Set-TextValue $ws.Range("E48") "  -0.20%  "
Set-TextValue $ws.Range("D49") "2.43"
Set-TextValue $ws.Range("E49") "  +5.79%  "
Set-TextValue $ws.Range("D50") "0.597"
Set-TextValue $ws.Range("E50") "  +1.19%  "
Set-TextValue $ws.Range("D51") "3.89"
Set-TextValue $ws.Range("E51") "  +0.06%  "
